$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a daily price log for "Albahaca" (Hortaliza) at
# "Femacal de La Calera", one row per observation, newest report inserted
# at the top of the data block (row 66, just below the header + first
# existing entry). Insert a new row at 66, which pushes the existing
# rows 66-172 down to 67-173 (row 173 ends up holding what used to be
# row 172's data), then populate the new row 66 with the new report.

$ws.Rows("66:66").Insert()

$ws.Range("A66").Value = 3
$ws.Range("B66").Value = "Femacal de La Calera"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value2 = 44665
$ws.Range("E66").Value = 5
$ws.Range("F66").Value = 100112052
$ws.Range("G66").Value = "Albahaca"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 125
$ws.Range("K66").Value = 4000
$ws.Range("L66").Value = 4500
$ws.Range("M66").Value = 4240
$ws.Range("N66").Value = "$/docena de matas"
$ws.Range("O66").Value = "Provincia de Quillota"
$ws.Range("P66").Value = 707
$ws.Range("Q66").Value = 6
$ws.Range("R66").Value = "Hortaliza"
